$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-27 22:18:14'
$ws.Range('H2').Value = "'56%"
$ws.Range('O2').Value = '5.4 °C'
$ws.Range('E3').Value = '2026-02-27 22:18:17'
$ws.Range('H3').Value = "'44%"
$ws.Range('O3').Value = '4.2 °C'
$ws.Range('E4').Value = '2026-02-27 22:18:19'
$ws.Range('E5').Value = '2026-02-27 22:18:22'
$ws.Range('H5').Value = "'46%"
$ws.Range('N5').Value = '0.4 °C 21:34 TU'
$ws.Range('O5').Value = '4.6 °C'
$ws.Range('E6').Value = '2026-02-27 22:18:24'
$ws.Range('O6').Value = '10.9 °C'
$ws.Range('E7').Value = '2026-02-27 22:18:26'
$ws.Range('E8').Value = '2026-02-27 22:18:29'
$ws.Range('O8').Value = '11.6 °C'
$ws.Range('E9').Value = '2026-02-27 22:18:31'
$ws.Range('H9').Value = "'91%"
$ws.Range('E10').Value = '2026-02-27 22:18:33'
$ws.Range('O10').Value = '10.8 °C'
$ws.Range('E11').Value = '2026-02-27 22:18:36'
$ws.Range('E12').Value = '2026-02-27 22:18:38'
$ws.Range('E13').Value = '2026-02-27 22:18:40'
$ws.Range('K13').Value = '15.0 MJ/m2'
$ws.Range('E14').Value = '2026-02-27 22:18:43'
$ws.Range('N14').Value = '5.4 °C 21:41 TU'
$ws.Range('O14').Value = '10.4 °C'
$ws.Range('E15').Value = '2026-02-27 22:18:45'
$ws.Range('O15').Value = '10.8 °C'
$ws.Range('E16').Value = '2026-02-27 22:18:47'
$ws.Range('N16').Value = '0.0 °C 21:38 TU'
$ws.Range('E17').Value = '2026-02-27 22:18:50'
$ws.Range('N17').Value = '4.8 °C 21:42 TU'
$ws.Range('O17').Value = '7.4 °C'
$ws.Range('E18').Value = '2026-02-27 22:18:52'
$ws.Range('H18').Value = "'82%"
$ws.Range('O18').Value = '11.9 °C'
$ws.Range('E19').Value = '2026-02-27 22:18:55'
$ws.Range('H19').Value = "'62%"
$ws.Range('E20').Value = '2026-02-27 22:18:57'
$ws.Range('E21').Value = '2026-02-27 22:18:59'
$ws.Range('O21').Value = '9.8 °C'
$ws.Range('E22').Value = '2026-02-27 22:19:02'
$ws.Range('K22').Value = '17.4 MJ/m2'
$ws.Range('E23').Value = '2026-02-27 22:19:04'
$ws.Range('H23').Value = "'43%"
$ws.Range('N23').Value = '0.7 °C 21:45 TU'
$ws.Range('E24').Value = '2026-02-27 22:19:07'
$ws.Range('J24').Value = '1023.3 hPa'
$ws.Range('E25').Value = '2026-02-27 22:19:09'
$ws.Range('N25').Value = '2.0 °C 21:50 TU'
$ws.Range('O25').Value = '5.9 °C'
$ws.Range('E26').Value = '2026-02-27 22:19:11'
$ws.Range('H26').Value = "'47%"
$ws.Range('N26').Value = '5.4 °C 21:59 TU'
$ws.Range('O26').Value = '10.0 °C'
$ws.Range('E27').Value = '2026-02-27 22:19:14'
$ws.Range('E28').Value = '2026-02-27 22:19:16'
$ws.Range('H28').Value = "'92%"
$ws.Range('O28').Value = '8.2 °C'
$ws.Range('E29').Value = '2026-02-27 22:19:19'
$ws.Range('E30').Value = '2026-02-27 22:19:21'
$ws.Range('O30').Value = '10.9 °C'
$ws.Range('E31').Value = '2026-02-27 22:19:23'
$ws.Range('E32').Value = '2026-02-27 22:19:25'
$ws.Range('H32').Value = "'55%"
$ws.Range('O32').Value = '7.8 °C'
$ws.Range('E33').Value = '2026-02-27 22:19:28'
$ws.Range('H33').Value = "'52%"
$ws.Range('O33').Value = '8.6 °C'
$ws.Range('E34').Value = '2026-02-27 22:19:30'
$ws.Range('O34').Value = '4.5 °C'
$ws.Range('E35').Value = '2026-02-27 22:19:33'
$ws.Range('H35').Value = "'42%"
$ws.Range('O35').Value = '11.9 °C'
$ws.Range('E36').Value = '2026-02-27 22:19:35'
$ws.Range('E37').Value = '2026-02-27 22:19:38'
$ws.Range('H37').Value = "'69%"
$ws.Range('E38').Value = '2026-02-27 22:19:40'
$ws.Range('E39').Value = '2026-02-27 22:19:42'
$ws.Range('H39').Value = "'32%"
$ws.Range('N39').Value = '1.2 °C 21:52 TU'
$ws.Range('O39').Value = '4.5 °C'
$ws.Range('E40').Value = '2026-02-27 22:19:45'
$ws.Range('H40').Value = "'68%"
$ws.Range('O40').Value = '8.9 °C'
$ws.Range('E41').Value = '2026-02-27 22:19:47'
$ws.Range('E42').Value = '2026-02-27 22:19:49'
$ws.Range('O42').Value = '11.3 °C'
$ws.Range('E43').Value = '2026-02-27 22:19:52'
$ws.Range('O43').Value = '9.3 °C'
$ws.Range('E44').Value = '2026-02-27 22:19:54'
$ws.Range('O44').Value = '1.9 °C'
$ws.Range('E45').Value = '2026-02-27 22:19:56'
$ws.Range('H45').Value = "'46%"
$ws.Range('O45').Value = '11.7 °C'
$ws.Range('E46').Value = '2026-02-27 22:19:59'
